# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3345
$wsExhibit.Range("F4").Value = 128
$wsExhibit.Range("F5").Value = 6943
$wsExhibit.Range("F6").Value = 2334
$wsExhibit.Range("F7").Value = 35
$wsExhibit.Range("F8").Value = 95
$wsExhibit.Range("F13").Value = 167
$wsExhibit.Range("F14").Value = 553

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3345
$wsAll.Range("F5").Value = 128
$wsAll.Range("F6").Value = 6943
$wsAll.Range("F7").Value = 2334
$wsAll.Range("F8").Value = 35
$wsAll.Range("F9").Value = 95
$wsAll.Range("F14").Value = 167
$wsAll.Range("F15").Value = 553
